$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3120081.25
$ws.Range("C9").Value = 490391.99
$ws.Range("D9").Value = 3610473.24
$ws.Range("E9").Value = 13.58248510380872
$ws.Range("F9").Value = 86.41751489619128
$ws.Range("G9").Value = -52.60605937369025
$ws.Range("H9").Value = -43.65564532490431
$ws.Range("I9").Value = 31188
$ws.Range("J9").Value = 1324
$ws.Range("K9").Value = 32512
$ws.Range("L9").Value = 22443
$ws.Range("M9").Value = 160.8730223232188
$ws.Range("N9").Value = 9.831288958249518
